# Updates cryptos list values to reflect the latest scrape
# (price and volume(1h) changes, plus two pairs of rows that swapped order)
# D-column cells are forced to Text format before assignment so that
# Excel's automatic numeric type inference doesn't strip trailing/leading
# zeros from the price strings (e.g. "4.600" -> 4.6).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.099.12'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.652.94'
$ws.Range("E3").Value = '  -1.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '218.71'
$ws.Range("E5").Value = '  -0.78%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5257'
$ws.Range("E6").Value = '  -1.19%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2673'
$ws.Range("E8").Value = '  +0.64%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06368'
$ws.Range("E9").Value = '  -0.08%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.52'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07686'
$ws.Range("E11").Value = '  -2.26%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.600'
$ws.Range("E12").Value = '  +1.29%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.660.05'
$ws.Range("E13").Value = '  -0.78%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.880.84'
$ws.Range("E14").Value = '  -1.05%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5615'
$ws.Range("E15").Value = '  -0.14%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₅8226'
$ws.Range("E16").Value = '  +1.12%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '65.57'
$ws.Range("E17").Value = '  -0.66%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '26.107.98'
$ws.Range("E18").Value = '  -1.00%  '
$ws.Range("E19").Value = '  -0.54%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.681'
$ws.Range("E20").Value = '  -1.08%  '
$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '190.84'
$ws.Range("E21").Value = '  -5.30%  '
$ws.Range("B22").Value = 'Avalanche'
$ws.Range("C22").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.31'
$ws.Range("E22").Value = '  -0.05%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.970'
$ws.Range("E23").Value = '  -1.66%  '
$ws.Range("E24").Value = '  -0.49%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.70'
$ws.Range("E25").Value = '  -0.78%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1202'
$ws.Range("E26").Value = '  -1.35%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.255'
$ws.Range("E27").Value = '  -0.13%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.96'
$ws.Range("E28").Value = '  -1.85%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.519'
$ws.Range("E29").Value = '  -0.01%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05641'
$ws.Range("E30").Value = '  -4.58%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.274'
$ws.Range("E31").Value = '  -1.07%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.492'
$ws.Range("E32").Value = '  -1.19%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.378'
$ws.Range("E33").Value = '  +1.35%  '
$ws.Range("E34").Value = '  -1.77%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.796'
$ws.Range("E35").Value = '  -1.29%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9464'
$ws.Range("E36").Value = '  -2.21%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.409'
$ws.Range("E37").Value = '  -0.91%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5778'
$ws.Range("E38").Value = '  -0.69%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01593'
$ws.Range("E39").Value = '  -1.73%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.965'
$ws.Range("E40").Value = '  -0.27%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8454'
$ws.Range("E41").Value = '  -1.82%  '
$ws.Range("E42").Value = '  -0.57%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.021.04'
$ws.Range("E43").Value = '  -5.36%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.24'
$ws.Range("E44").Value = '  -1.91%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.791.54'
$ws.Range("E45").Value = '  -1.05%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '58.35'
$ws.Range("E46").Value = '  -0.48%  '
$ws.Range("B47").Value = 'Frax'
$ws.Range("C47").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.003'
$ws.Range("E47").Value = '  -1.37%  '
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.05330'
$ws.Range("E48").Value = '  +3.50%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.047'
$ws.Range("E49").Value = '  -0.56%  '
$ws.Range("B50").Value = 'Mantle'
$ws.Range("C50").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4343'
$ws.Range("E50").Value = '  -1.71%  '
$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0₈102'
$ws.Range("E51").Value = '  -3.71%  '
